$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")
$ws.Activate()

# Rename scenario: "ClickOnCloseGroup" -> "ClickOnDeleteButton"
$ws.Range("A6").Value = "ClickOnDeleteButton"

# Fix group name typo: "Group4097" -> "Grou4097"
$ws.Range("B2").Value = "Grou4097"

# Update current selection to B3
$ws.Range("B3").Select()
